$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: update title (column D) for "ai plus(est soft)" feed
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 37: update title (column D) and link (column E) for dsba_seminar feed
$ws.Range("D37").Value = "[Paper Review] Structure Extraction in Task-Oriented Dialogues with Slot Clustering"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1959&mod=document&pageid=1"

# Row 46: update title (column D) and link (column E) for BioinformaticsAndMe feed
$ws.Range("D46").Value = "Child-Pugh Score (차일드-퍼 점수)"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/440"
